$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product name text (shared string) on both sheets (B1)
$newProductName = "2596-RBI-EI-DB-DL-REC-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-OVERDUE-FEE-%LOANAMT+INT-MORE-AMT-1st"
$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Update shortname (B2) from numeric 2596 to text "259w"
$wsInput.Range("B2").Value = "259w"

# Update the selection/view on the input sheet: select B2, and make the
# input sheet the active/selected one (matches tabSelected in target)
$wsInput.Activate()
$wsInput.Range("B2").Select()
